# "Support direct range notation" - adds a new test case row to the
# `model` sheet that uses an unwrapped (direct) multi-area range reference
# as a formula, and adds a "winter" data row + a duplicated data table on
# the `charts` sheet for that new test case to reference. Named ranges
# that describe the charts table are widened by one row to match.

$wb = $excel.ActiveWorkbook
$wsModel = $wb.Worksheets.Item("model")
$wsCharts = $wb.Worksheets.Item("charts")

# --- charts sheet: add the "winter" row + a duplicated table for the -------
# --- new direct-range-notation formula to reference -------------------------
$wsCharts.Range("A5").Value = "冬"
$wsCharts.Range("B5").Value = 130
$wsCharts.Range("C5").Value = 0
$wsCharts.Range("D5").Value = 0.6

$wsCharts.Range("A7").Value = "Season"
$wsCharts.Range("B7").Value = "売り上げ"
$wsCharts.Range("C7").Value = "利益"
$wsCharts.Range("D7").Value = "利益率"

$wsCharts.Range("A8").Value = "春"
$wsCharts.Range("B8").Value = 100
$wsCharts.Range("C8").Value = 50
$wsCharts.Range("D8").Value = 0.5

$wsCharts.Range("A9").Value = "夏"
$wsCharts.Range("B9").Value = 110
$wsCharts.Range("C9").Value = 60
$wsCharts.Range("D9").Value = 0.5

$wsCharts.Range("A10").Value = "秋"
$wsCharts.Range("B10").Value = 120
$wsCharts.Range("C10").Value = 70
$wsCharts.Range("D10").Value = 0.5

$wsCharts.Range("A11").Value = "冬"
$wsCharts.Range("B11").Value = 130
$wsCharts.Range("C11").Value = 0
$wsCharts.Range("D11").Value = 0.6

# --- model sheet: insert a new "normal2" test row at row 4 -----------------
$wsModel.Rows.Item(4).Insert()
$wsModel.Range("A4").Value = "p02"
$wsModel.Range("B4").Value = "normal2"
$wsModel.Range("D4").Formula = "=charts!A7:D9,charts!A10:D11"

# --- defined names: widen the charts ranges by the new winter row ----------
$wb.Names.Item("chart_sidebyside").RefersTo = "=charts!`$A`$1:`$A`$5,charts!`$D`$1:`$D`$5"
$wb.Names.Item("chart01").RefersTo = "=charts!`$A`$1:`$D`$5"

# --- restore the selections shown in the target workbook -------------------
$wsModel.Activate() | Out-Null
$wsModel.Range("D5").Select() | Out-Null

$wsCharts.Activate() | Out-Null
$wsCharts.Range("E13").Select() | Out-Null

$wsModel.Activate() | Out-Null
